# Trade #112 closed at 2026-02-17 09:19:37 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Summary sheet updates
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.16   # Current Capital
$summary.Range("B4").Value = 0.17      # Total P&L $
$summary.Range("B6").Value = 112       # Total Trades
$summary.Range("B7").Value = 49        # Winning Trades
$summary.Range("B9").Value = 43.75     # Win Rate %

# ---------------------------------------------------------------
# 2. Strategy Status sheet updates (MarketMaking is row 4)
# ---------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.16     # Capital
$status.Range("D4").Value = 112        # Trades
$status.Range("E4").Value = 0.17       # P&L $
$status.Range("F4").Value = 0.16       # P&L %
$status.Range("G4").Value = 43.75      # Win Rate %

# ---------------------------------------------------------------
# 3. New trade row (#112) appended as row 113 to both the
#    "All Trades" log and the per-strategy "MarketMaking" log.
# ---------------------------------------------------------------
foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item(113, 1).Value = 112                 # Trade #

    # Date column looks like "2026-02-17" - without the leading
    # apostrophe Excel auto-detects it as a date serial, which the
    # source file never does (all dates/times are kept as plain text).
    $ws.Cells.Item(113, 2).Value = "'2026-02-17"
    $ws.Cells.Item(113, 2).ClearFormats()

    $ws.Cells.Item(113, 3).Value = "09:19:30"          # Time
    $ws.Cells.Item(113, 4).Value = "MarketMaking"      # Strategy
    $ws.Cells.Item(113, 5).Value = "UP"                # Side
    $ws.Cells.Item(113, 6).Value = 0.22                # Entry Price
    $ws.Cells.Item(113, 7).Value = 0.23                # Exit Price
    $ws.Cells.Item(113, 8).Value = "CLOSED"            # Status
    $ws.Cells.Item(113, 9).Value = 4.5455              # P&L %
    $ws.Cells.Item(113, 10).Value = 0.01               # P&L $
    $ws.Cells.Item(113, 11).Value = 100.16             # Capital After
    $ws.Cells.Item(113, 12).Value = 0                  # Entry Slippage (bps)
    $ws.Cells.Item(113, 13).Value = 0                  # Exit Slippage (bps)
    $ws.Cells.Item(113, 14).Value = 0.6                # Confidence
    $ws.Cells.Item(113, 15).Value = "Normal spread capture: 19600 bps"  # Entry Reason
    $ws.Cells.Item(113, 16).Value = "early_exit"       # Exit Reason
    $ws.Cells.Item(113, 17).Value = 0.14               # Duration (min)
}
